# Auto-generated script to update Leve profit calculation columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR Table_* sheets, reflecting
# refreshed market-board pricing data from the scheduled runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 11
$ws.Cells.Item(11, 8).Value = 519.7646999999999   # H11: 493.16666 -> 519.7646999999999
$ws.Cells.Item(11, 9).Value = 519.7646999999999   # I11: 493.16666 -> 519.7646999999999
$ws.Cells.Item(11, 11).Value = 519.7646999999999   # K11: 493.16666 -> 519.7646999999999
$ws.Cells.Item(11, 13).Value = -379.7646999999999   # M11: -353.16666 -> -379.7646999999999

# ALC row 49
$ws.Cells.Item(49, 8).Value = 517   # H49: 2469.6 -> 517
$ws.Cells.Item(49, 9).Value = 517   # I49: 1922 -> 517
$ws.Cells.Item(49, 10).Value = 0   # J49: 3291 -> 0
$ws.Cells.Item(49, 11).Value = 1551   # K49: 5766 -> 1551
$ws.Cells.Item(49, 12).Value = 0   # L49: 9873 -> 0
$ws.Cells.Item(49, 13).Value = -1415   # M49: -5630 -> -1415
$ws.Cells.Item(49, 14).ClearContents()   # N49: -10145 -> (removed)

# ALC row 70
$ws.Cells.Item(70, 8).Value = 32160716   # H70: 33135252 -> 32160716
$ws.Cells.Item(70, 9).Value = 13004837   # I70: 13689250 -> 13004837
$ws.Cells.Item(70, 11).Value = 39014511   # K70: 41067750 -> 39014511
$ws.Cells.Item(70, 13).Value = -39014241   # M70: -41067480 -> -39014241

# ALC row 73
$ws.Cells.Item(73, 8).Value = 32160716   # H73: 33135252 -> 32160716
$ws.Cells.Item(73, 9).Value = 13004837   # I73: 13689250 -> 13004837
$ws.Cells.Item(73, 11).Value = 39014511   # K73: 41067750 -> 39014511
$ws.Cells.Item(73, 13).Value = -39013575   # M73: -41066814 -> -39013575

# ALC row 116
$ws.Cells.Item(116, 8).Value = 25009180   # H116: 19239292 -> 25009180
$ws.Cells.Item(116, 9).Value = 250000000   # I116: 83336340 -> 250000000
$ws.Cells.Item(116, 10).Value = 10200   # J116: 10179.9 -> 10200
$ws.Cells.Item(116, 11).Value = 250000000   # K116: 83336340 -> 250000000
$ws.Cells.Item(116, 12).Value = 10200   # L116: 10179.9 -> 10200
$ws.Cells.Item(116, 13).Value = -249996558   # M116: -83332898 -> -249996558
$ws.Cells.Item(116, 14).Value = -17084   # N116: -17063.9 -> -17084

# ALC row 118
$ws.Cells.Item(118, 8).Value = 4591.6   # H118: 2083.3333 -> 4591.6
$ws.Cells.Item(118, 9).Value = 916.3333   # I118: 999 -> 916.3333
$ws.Cells.Item(118, 10).Value = 10104.5   # J118: 2300.2 -> 10104.5
$ws.Cells.Item(118, 11).Value = 2748.9999   # K118: 2997 -> 2748.9999
$ws.Cells.Item(118, 12).Value = 30313.5   # L118: 6900.599999999999 -> 30313.5
$ws.Cells.Item(118, 13).Value = -1091.9999   # M118: -1340 -> -1091.9999
$ws.Cells.Item(118, 14).Value = -33627.5   # N118: -10214.6 -> -33627.5

# ALC row 135
$ws.Cells.Item(135, 8).Value = 435892.56   # H135: 455674.4 -> 435892.56
$ws.Cells.Item(135, 9).Value = 501093.44   # I135: 527430.4 -> 501093.44
$ws.Cells.Item(135, 11).Value = 4509840.96   # K135: 4746873.600000001 -> 4509840.96
$ws.Cells.Item(135, 13).Value = -4507305.96   # M135: -4744338.600000001 -> -4507305.96

# ALC row 137
$ws.Cells.Item(137, 8).Value = 3805.6   # H137: 2892.9768 -> 3805.6
$ws.Cells.Item(137, 9).Value = 5043   # I137: 2263.3809 -> 5043
$ws.Cells.Item(137, 10).Value = 3429   # J137: 3493.9546 -> 3429
$ws.Cells.Item(137, 11).Value = 15129   # K137: 6790.1427 -> 15129
$ws.Cells.Item(137, 12).Value = 10287   # L137: 10481.8638 -> 10287
$ws.Cells.Item(137, 13).Value = -12579   # M137: -4240.1427 -> -12579
$ws.Cells.Item(137, 14).Value = -15387   # N137: -15581.8638 -> -15387

# ALC row 138
$ws.Cells.Item(138, 8).Value = 1643323.6   # H138: 1670680.6 -> 1643323.6
$ws.Cells.Item(138, 10).Value = 2637512.8   # J138: 2708745.5 -> 2637512.8
$ws.Cells.Item(138, 12).Value = 7912538.399999999   # L138: 8126236.5 -> 7912538.399999999
$ws.Cells.Item(138, 14).Value = -7922818.399999999   # N138: -8136516.5 -> -7922818.399999999

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Cells.Item(32, 8).Value = 2322416.2   # H32: 2322429.8 -> 2322416.2
$ws.Cells.Item(32, 9).Value = 2722731.8   # I32: 2722747.5 -> 2722731.8
$ws.Cells.Item(32, 11).Value = 2722731.8   # K32: 2722747.5 -> 2722731.8
$ws.Cells.Item(32, 13).Value = -2722444.8   # M32: -2722460.5 -> -2722444.8

# ARM row 74
$ws.Cells.Item(74, 8).Value = 27251.08   # H74: 28345 -> 27251.08
$ws.Cells.Item(74, 9).Value = 35982.055   # I74: 38040 -> 35982.055
$ws.Cells.Item(74, 11).Value = 35982.055   # K74: 38040 -> 35982.055
$ws.Cells.Item(74, 13).Value = -35108.055   # M74: -37166 -> -35108.055

# ARM row 77
$ws.Cells.Item(77, 8).Value = 27251.08   # H77: 28345 -> 27251.08
$ws.Cells.Item(77, 9).Value = 35982.055   # I77: 38040 -> 35982.055
$ws.Cells.Item(77, 11).Value = 179910.275   # K77: 190200 -> 179910.275
$ws.Cells.Item(77, 13).Value = -175542.275   # M77: -185832 -> -175542.275

# ARM row 132
$ws.Cells.Item(132, 8).Value = 3288.8667   # H132: 4142.1885 -> 3288.8667
$ws.Cells.Item(132, 9).Value = 1932.898   # I132: 1999.8043 -> 1932.898
$ws.Cells.Item(132, 10).Value = 9329.091   # J132: 8426.956 -> 9329.091
$ws.Cells.Item(132, 11).Value = 5798.694   # K132: 5999.4129 -> 5798.694
$ws.Cells.Item(132, 12).Value = 27987.273   # L132: 25280.868 -> 27987.273
$ws.Cells.Item(132, 13).Value = -3268.694   # M132: -3469.4129 -> -3268.694
$ws.Cells.Item(132, 14).Value = -33047.273   # N132: -30340.868 -> -33047.273

$ws = $wb.Worksheets.Item("BSM")
# BSM row 26
$ws.Cells.Item(26, 8).Value = 32133   # H26: 36066.6 -> 32133
$ws.Cells.Item(26, 9).Value = 21228.5   # I26: 24149.666 -> 21228.5
$ws.Cells.Item(26, 11).Value = 21228.5   # K26: 24149.666 -> 21228.5
$ws.Cells.Item(26, 13).Value = -20936.5   # M26: -23857.666 -> -20936.5

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Cells.Item(31, 8).Value = 5959047   # H31: 6296244 -> 5959047
$ws.Cells.Item(31, 9).Value = 2905.0667   # I31: 3031.9644 -> 2905.0667
$ws.Cells.Item(31, 10).Value = 12831519   # J31: 13344642 -> 12831519
$ws.Cells.Item(31, 11).Value = 2905.0667   # K31: 3031.9644 -> 2905.0667
$ws.Cells.Item(31, 12).Value = 12831519   # L31: 13344642 -> 12831519
$ws.Cells.Item(31, 13).Value = -2610.0667   # M31: -2736.9644 -> -2610.0667
$ws.Cells.Item(31, 14).Value = -12832109   # N31: -13345232 -> -12832109

# CRP row 34
$ws.Cells.Item(34, 8).Value = 5959047   # H34: 6296244 -> 5959047
$ws.Cells.Item(34, 9).Value = 2905.0667   # I34: 3031.9644 -> 2905.0667
$ws.Cells.Item(34, 10).Value = 12831519   # J34: 13344642 -> 12831519
$ws.Cells.Item(34, 11).Value = 2905.0667   # K34: 3031.9644 -> 2905.0667
$ws.Cells.Item(34, 12).Value = 12831519   # L34: 13344642 -> 12831519
$ws.Cells.Item(34, 13).Value = -2703.0667   # M34: -2829.9644 -> -2703.0667
$ws.Cells.Item(34, 14).Value = -12831923   # N34: -13345046 -> -12831923

# CRP row 94
$ws.Cells.Item(94, 8).Value = 1208.25   # H94: 1211.3 -> 1208.25
$ws.Cells.Item(94, 9).Value = 1660.1666   # I94: 1936.4 -> 1660.1666
$ws.Cells.Item(94, 10).Value = 1014.5714   # J94: 969.6 -> 1014.5714
$ws.Cells.Item(94, 11).Value = 1660.1666   # K94: 1936.4 -> 1660.1666
$ws.Cells.Item(94, 12).Value = 1014.5714   # L94: 969.6 -> 1014.5714
$ws.Cells.Item(94, 13).Value = -1209.1666   # M94: -1485.4 -> -1209.1666
$ws.Cells.Item(94, 14).Value = -1916.5714   # N94: -1871.6 -> -1916.5714

# CRP row 99
$ws.Cells.Item(99, 8).Value = 8366.556   # H99: 8357.277 -> 8366.556
$ws.Cells.Item(99, 9).Value = 11116.167   # I99: 12639.4 -> 11116.167
$ws.Cells.Item(99, 10).Value = 6991.75   # J99: 6710.3076 -> 6991.75
$ws.Cells.Item(99, 11).Value = 11116.167   # K99: 12639.4 -> 11116.167
$ws.Cells.Item(99, 12).Value = 6991.75   # L99: 6710.3076 -> 6991.75
$ws.Cells.Item(99, 13).Value = -9618.166999999999   # M99: -11141.4 -> -9618.166999999999
$ws.Cells.Item(99, 14).Value = -9987.75   # N99: -9706.3076 -> -9987.75

# CRP row 126
$ws.Cells.Item(126, 8).Value = 8366.556   # H126: 8357.277 -> 8366.556
$ws.Cells.Item(126, 9).Value = 11116.167   # I126: 12639.4 -> 11116.167
$ws.Cells.Item(126, 10).Value = 6991.75   # J126: 6710.3076 -> 6991.75
$ws.Cells.Item(126, 11).Value = 33348.501   # K126: 37918.2 -> 33348.501
$ws.Cells.Item(126, 12).Value = 20975.25   # L126: 20130.9228 -> 20975.25
$ws.Cells.Item(126, 13).Value = -30878.501   # M126: -35448.2 -> -30878.501
$ws.Cells.Item(126, 14).Value = -25915.25   # N126: -25070.9228 -> -25915.25

# CRP row 132
$ws.Cells.Item(132, 8).Value = 6456552   # H132: 6456553 -> 6456552
$ws.Cells.Item(132, 9).Value = 2894.0908   # I132: 2896.2424 -> 2894.0908
$ws.Cells.Item(132, 11).Value = 8682.2724   # K132: 8688.727200000001 -> 8682.2724
$ws.Cells.Item(132, 13).Value = -6152.2724   # M132: -6158.727200000001 -> -6152.2724

# CRP row 134
$ws.Cells.Item(134, 8).Value = 5841.039   # H134: 5939 -> 5841.039
$ws.Cells.Item(134, 9).Value = 2659.2173   # I134: 2737.2273 -> 2659.2173
$ws.Cells.Item(134, 11).Value = 7977.651899999999   # K134: 8211.6819 -> 7977.651899999999
$ws.Cells.Item(134, 13).Value = -5442.651899999999   # M134: -5676.6819 -> -5442.651899999999

# CRP row 141
$ws.Cells.Item(141, 8).Value = 137995.6   # H141: 181999.33 -> 137995.6
$ws.Cells.Item(141, 10).Value = 137995.6   # J141: 181999.33 -> 137995.6
$ws.Cells.Item(141, 12).Value = 137995.6   # L141: 181999.33 -> 137995.6
$ws.Cells.Item(141, 14).Value = -148355.6   # N141: -192359.33 -> -148355.6

$ws = $wb.Worksheets.Item("CUL")
# CUL row 4
$ws.Cells.Item(4, 8).Value = 87120504   # H4: 93822060 -> 87120504
$ws.Cells.Item(4, 9).Value = 121100120   # I4: 134555650 -> 121100120
$ws.Cells.Item(4, 11).Value = 363300360   # K4: 403666950 -> 363300360
$ws.Cells.Item(4, 13).Value = -363300248   # M4: -403666838 -> -363300248

# CUL row 26
$ws.Cells.Item(26, 8).Value = 303.8   # H26: 299.70834 -> 303.8
$ws.Cells.Item(26, 9).Value = 179.8   # I26: 159.2 -> 179.8
$ws.Cells.Item(26, 10).Value = 334.8   # J26: 336.6842 -> 334.8
$ws.Cells.Item(26, 11).Value = 539.4000000000001   # K26: 477.6 -> 539.4000000000001
$ws.Cells.Item(26, 12).Value = 1004.4   # L26: 1010.0526 -> 1004.4
$ws.Cells.Item(26, 13).Value = -251.4000000000001   # M26: -189.6 -> -251.4000000000001
$ws.Cells.Item(26, 14).Value = -1580.4   # N26: -1586.0526 -> -1580.4

# CUL row 123
$ws.Cells.Item(123, 8).Value = 3380   # H123: 4120 -> 3380
$ws.Cells.Item(123, 10).Value = 7650   # J123: 9500 -> 7650
$ws.Cells.Item(123, 12).Value = 22950   # L123: 28500 -> 22950
$ws.Cells.Item(123, 14).Value = -27850   # N123: -33400 -> -27850

# CUL row 132
$ws.Cells.Item(132, 8).Value = 10479.414   # H132: 10507 -> 10479.414
$ws.Cells.Item(132, 9).Value = 4000   # I132: 4164.706 -> 4000
$ws.Cells.Item(132, 10).Value = 21082.092   # J132: 19491.916 -> 21082.092
$ws.Cells.Item(132, 11).Value = 36000   # K132: 37482.354 -> 36000
$ws.Cells.Item(132, 12).Value = 189738.828   # L132: 175427.244 -> 189738.828
$ws.Cells.Item(132, 13).Value = -33470   # M132: -34952.354 -> -33470
$ws.Cells.Item(132, 14).Value = -194798.828   # N132: -180487.244 -> -194798.828

# CUL row 139
$ws.Cells.Item(139, 8).Value = 4612.472   # H139: 4579.8423 -> 4612.472
$ws.Cells.Item(139, 9).Value = 2523.913   # I139: 2459.7917 -> 2523.913
$ws.Cells.Item(139, 10).Value = 8307.615   # J139: 8214.214 -> 8307.615
$ws.Cells.Item(139, 11).Value = 7571.739   # K139: 7379.375100000001 -> 7571.739
$ws.Cells.Item(139, 12).Value = 24922.845   # L139: 24642.642 -> 24922.845
$ws.Cells.Item(139, 13).Value = -2431.739   # M139: -2239.375100000001 -> -2431.739
$ws.Cells.Item(139, 14).Value = -35202.845   # N139: -34922.642 -> -35202.845

$ws = $wb.Worksheets.Item("GSM")
# GSM row 113
$ws.Cells.Item(113, 8).Value = 5951.978   # H113: 6158.386 -> 5951.978
$ws.Cells.Item(113, 9).Value = 2718.1428   # I113: 2855.7368 -> 2718.1428
$ws.Cells.Item(113, 11).Value = 2718.1428   # K113: 2855.7368 -> 2718.1428
$ws.Cells.Item(113, 13).Value = -548.1428000000001   # M113: -685.7368000000001 -> -548.1428000000001

# GSM row 126
$ws.Cells.Item(126, 8).Value = 6286.8945   # H126: 6444.7896 -> 6286.8945
$ws.Cells.Item(126, 9).Value = 5904.077   # I126: 6134.846 -> 5904.077
$ws.Cells.Item(126, 11).Value = 17712.231   # K126: 18404.538 -> 17712.231
$ws.Cells.Item(126, 13).Value = -15242.231   # M126: -15934.538 -> -15242.231

# GSM row 132
$ws.Cells.Item(132, 8).Value = 5055.425   # H132: 4556.8936 -> 5055.425
$ws.Cells.Item(132, 9).Value = 2085.3845   # I132: 2005.3636 -> 2085.3845
$ws.Cells.Item(132, 11).Value = 6256.1535   # K132: 6016.0908 -> 6256.1535
$ws.Cells.Item(132, 13).Value = -3726.1535   # M132: -3486.0908 -> -3726.1535

# GSM row 135
$ws.Cells.Item(135, 8).Value = 51771.25   # H135: 58030.91 -> 51771.25
$ws.Cells.Item(135, 10).Value = 51771.25   # J135: 58030.91 -> 51771.25
$ws.Cells.Item(135, 12).Value = 51771.25   # L135: 58030.91 -> 51771.25
$ws.Cells.Item(135, 14).Value = -61911.25   # N135: -68170.91 -> -61911.25

$ws = $wb.Worksheets.Item("LTW")
# LTW row 132
$ws.Cells.Item(132, 8).Value = 7818150.5   # H132: 7468151.5 -> 7818150.5
$ws.Cells.Item(132, 9).Value = 15154033   # I132: 13891322 -> 15154033
$ws.Cells.Item(132, 11).Value = 45462099   # K132: 41673966 -> 45462099
$ws.Cells.Item(132, 13).Value = -45459569   # M132: -41671436 -> -45459569

$ws = $wb.Worksheets.Item("WVR")
# WVR row 56
$ws.Cells.Item(56, 8).Value = 69000   # H56: 66666 -> 69000
$ws.Cells.Item(56, 9).Value = 69000   # I56: 0 -> 69000
$ws.Cells.Item(56, 10).Value = 0   # J56: 66666 -> 0
$ws.Cells.Item(56, 11).Value = 69000   # K56: 0 -> 69000
$ws.Cells.Item(56, 12).Value = 0   # L56: 66666 -> 0
$ws.Cells.Item(56, 13).Value = -68286   # M56: (new) -> -68286
$ws.Cells.Item(56, 14).ClearContents()   # N56: -68094 -> (removed)

# WVR row 113
$ws.Cells.Item(113, 8).Value = 8989.200000000001   # H113: 9300.344999999999 -> 8989.200000000001
$ws.Cells.Item(113, 9).Value = 12607.81   # I113: 12052.909 -> 12607.81
$ws.Cells.Item(113, 10).Value = 545.7778   # J113: 649.4286 -> 545.7778
$ws.Cells.Item(113, 11).Value = 37823.43   # K113: 36158.727 -> 37823.43
$ws.Cells.Item(113, 12).Value = 1637.3334   # L113: 1948.2858 -> 1637.3334
$ws.Cells.Item(113, 13).Value = -35653.43   # M113: -33988.727 -> -35653.43
$ws.Cells.Item(113, 14).Value = -5977.3334   # N113: -6288.2858 -> -5977.3334

# WVR row 132
$ws.Cells.Item(132, 8).Value = 29435100   # H132: 23828508 -> 29435100
$ws.Cells.Item(132, 9).Value = 45465480   # I132: 33341486 -> 45465480
$ws.Cells.Item(132, 11).Value = 136396440   # K132: 100024458 -> 136396440
$ws.Cells.Item(132, 13).Value = -136393910   # M132: -100021928 -> -136393910

# WVR row 136
$ws.Cells.Item(136, 8).Value = 23837458   # H136: 23837564 -> 23837458
$ws.Cells.Item(136, 9).Value = 43479136   # I136: 45455436 -> 43479136
$ws.Cells.Item(136, 10).Value = 60692.316   # J136: 57907.7 -> 60692.316
$ws.Cells.Item(136, 11).Value = 130437408   # K136: 136366308 -> 130437408
$ws.Cells.Item(136, 12).Value = 182076.948   # L136: 173723.1 -> 182076.948
$ws.Cells.Item(136, 13).Value = -130434858   # M136: -136363758 -> -130434858
$ws.Cells.Item(136, 14).Value = -187176.948   # N136: -178823.1 -> -187176.948
